# Fixing network data cleaning scripts
# - Rename header columns to machine-friendly snake_case names
# - Normalize "de"/"del" -> "De"/"Del" (and "la" -> "La") in a handful of
#   municipality / state names
# - Tweak a floating point rounding artifact in D87
# - Drop the trailing sample-size / source / footer rows (101-105, 476-480)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames (row 1) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Capitalization fixes ---
$ws.Range("B8").Value  = "Amatenango De La Frontera"
$ws.Range("A17").Value = "Ciudad De México"
$ws.Range("A28").Value = "Estado De México"
$ws.Range("B29").Value = "Atizapán De Zaragoza"
$ws.Range("B35").Value = "Jaral Del Progreso"
$ws.Range("B43").Value = "Cuautepec De Hinojosa"
$ws.Range("B44").Value = "Tulancingo De Bravo"
$ws.Range("B46").Value = "Autlán De Navarro"
$ws.Range("B51").Value = "Zapotlán Del Rey"
$ws.Range("B55").Value = "Tiquicheo De Nicolás Romero"
$ws.Range("B67").Value = "Teotitlán De Flores Magón"
$ws.Range("B93").Value = "Vega De Alatorre"

# --- Minor floating point value tweak ---
$ws.Range("D87").Value = 0.09090909090909093

# --- Remove trailer rows (delete the higher block first so the lower
#     block's row numbers are unaffected by the shift) ---
$ws.Range("A476:A480").EntireRow.Delete()
$ws.Range("A101:A105").EntireRow.Delete()
